# Updates cryptos list price/volume figures (GitHub Actions scheduled refresh).
# Source diff only changes column D (Price) and column E (Volume(1h)) text values
# for rows 2-51; all other cells/columns are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "29.216.94"
$ws.Range("E2").Value = "  -0.55%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.859.70"
$ws.Range("E3").Value = "  -0.66%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5: XRP
$ws.Range("E5").Value = "  -0.36%  "

# Row 6: BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.42"
$ws.Range("E6").Value = "  +0.60%  "

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8: Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07763"
$ws.Range("E8").Value = "  -0.59%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3076"
$ws.Range("E9").Value = "  +0.31%  "

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.12"
$ws.Range("E10").Value = "  -0.20%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08267"
$ws.Range("E11").Value = "  +0.16%  "

# Row 12: WrappedEther
$ws.Range("D12").Value = "1.851.96"
$ws.Range("E12").Value = "  -0.89%  "

# Row 13: Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.229"
$ws.Range("E13").Value = "  -0.18%  "

# Row 14: Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7162"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15: Litecoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.38"
$ws.Range("E15").Value = "  -0.17%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "29.190.26"
$ws.Range("E16").Value = "  -0.70%  "

# Row 17: Uniswap
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.862"
$ws.Range("E17").Value = "  +0.24%  "

# Row 18: BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "244.06"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007795"
$ws.Range("E19").Value = "  -0.90%  "

# Row 20: Avalanche
$ws.Range("E20").Value = "  -0.91%  "

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.103.27"
$ws.Range("E21").Value = "  -0.43%  "

# Row 22: Dai
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"

# Row 23: Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.966"
$ws.Range("E23").Value = "  +2.55%  "

# Row 24: BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.06%  "

# Row 25: Stellar
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1590"
$ws.Range("E25").Value = "  +2.77%  "

# Row 26: Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.55"
$ws.Range("E26").Value = "  -0.65%  "

# Row 27: Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.908"
$ws.Range("E27").Value = "  -0.97%  "

# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.33"
$ws.Range("E28").Value = "  +0.12%  "

# Row 29: PancakeSwap
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.496"
$ws.Range("E29").Value = "  +0.95%  "

# Row 30: Toncoin
$ws.Range("E30").Value = "  -3.12%  "

# Row 31: Filecoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.416"
$ws.Range("E31").Value = "  +2.08%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.228"
$ws.Range("E32").Value = "  +3.43%  "

# Row 33: Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05187"
$ws.Range("E33").Value = "  -1.16%  "

# Row 34: LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.911"
$ws.Range("E34").Value = "  -0.95%  "

# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.172"
$ws.Range("E35").Value = "  -2.08%  "

# Row 36: ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7266"
$ws.Range("E36").Value = "  +1.48%  "

# Row 37: HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("E37").Value = "  -0.23%  "

# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01853"
$ws.Range("E38").Value = "  -0.69%  "

# Row 39: MXToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.690"
$ws.Range("E39").Value = "  -1.15%  "

# Row 40: Maker
$ws.Range("D40").Value = "1.164.63"
$ws.Range("E40").Value = "  -1.32%  "

# Row 41: TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9061"
$ws.Range("E41").Value = "  -0.30%  "

# Row 42: FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.149"
$ws.Range("E42").Value = "  +2.23%  "

# Row 43: Aave
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.22"
$ws.Range("E43").Value = "  -0.13%  "

# Row 44: PaxDollar
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.09%  "

# Row 45: Quant
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.59"
$ws.Range("E45").Value = "  -0.80%  "

# Row 46: RocketPoolETH
$ws.Range("D46").Value = "2.001.05"
$ws.Range("E46").Value = "  -0.53%  "

# Row 47: Mantle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5220"
$ws.Range("E47").Value = "  -2.77%  "

# Row 48: RenderToken
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.766"
$ws.Range("E48").Value = "  +0.16%  "

# Row 49: BabyDogeCoin
$ws.Range("E49").Value = "  -1.33%  "

# Row 50: EnergySwap
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.310"
$ws.Range("E50").Value = "  +1.79%  "

# Row 51: SynthetixNetwork
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.858"
$ws.Range("E51").Value = "  +1.08%  "
